$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = 9871.99
$ws.Range("B9").Value = 9815.06
$ws.Range("C9").Value = 307.21
$ws.Range("D9").Value = 305.42
$ws.Range("E9").Value = $true
$ws.Range("F9").Value = -0.58
$ws.Range("G9").Value = 42609.488981481481
$ws.Range("G9").NumberFormat = "m/d/yy h:mm"
$ws.Range("H9").Value = $true
